# "Added checks and writing results to excel"
# The Update/Delete automated checks for the first 9 Zone test rows (2-10)
# came back failing, so the results are recorded as "Fail" instead of the
# previous placeholder FALSE values. The remaining Category test rows
# (11-24) had their Create/Read checks re-run and now pass, so those cells
# flip from FALSE to TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 4).Value = "Fail"
    $ws.Cells.Item($r, 5).Value = "Fail"
}

for ($r = 11; $r -le 24; $r++) {
    $ws.Cells.Item($r, 2).Value = $true
    $ws.Cells.Item($r, 3).Value = $true
}

$ws.Range("D2").Select()
